$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10 (ALC)
$ws.Range("H10").Value = 300
$ws.Range("I10").Value = 300
$ws.Range("K10").Value = 300
$ws.Range("M10").Value = -7

# Row 17 (ALC)
$ws.Range("H17").Value = 184172.11
$ws.Range("J17").Value = 189405.6
$ws.Range("L17").Value = 568216.8
$ws.Range("N17").Value = -568552.8

# Row 92 (ALC)
$ws.Range("H92").Value = 1465906
$ws.Range("I92").Value = 679827.9
$ws.Range("K92").Value = 679827.9
$ws.Range("M92").Value = -678579.9

# Row 97 (ALC)
$ws.Range("H97").Value = 5633.6
$ws.Range("J97").Value = 5633.6
$ws.Range("L97").Value = 16900.8
$ws.Range("N97").Value = -17892.8

# Row 99 (ALC)
$ws.Range("H99").Value = 200000220
$ws.Range("I99").Value = 272.25
$ws.Range("K99").Value = 816.75
$ws.Range("M99").Value = 681.25

# Row 101 (ALC)
$ws.Range("H101").Value = 3999.75
$ws.Range("I101").Value = 3000
$ws.Range("J101").Value = 4333
$ws.Range("K101").Value = 9000
$ws.Range("L101").Value = 12999
$ws.Range("M101").Value = -7378
$ws.Range("N101").Value = -16243

# Row 106 (ALC)
$ws.Range("H106").Value = 66670452
$ws.Range("I106").Value = 71431200
$ws.Range("K106").Value = 71431200
$ws.Range("M106").Value = -71430569

$ws = $wb.Worksheets.Item("ARM")
# Row 80 (ARM)
$ws.Range("H80").Value = 17999.6
$ws.Range("J80").Value = 49998
$ws.Range("L80").Value = 49998
$ws.Range("N80").Value = -51994

# Row 83 (ARM)
$ws.Range("H83").Value = 17999.6
$ws.Range("J83").Value = 49998
$ws.Range("L83").Value = 149994
$ws.Range("N83").Value = -159978

# Row 102 (ARM)
$ws.Range("H102").Value = 3173
$ws.Range("I102").Value = 3173
$ws.Range("K102").Value = 3173
$ws.Range("M102").Value = -1551

$ws = $wb.Worksheets.Item("BSM")
# Row 40 (BSM)
$ws.Range("H40").Value = 35000
$ws.Range("J40").Value = 35000
$ws.Range("L40").Value = 35000
$ws.Range("N40").Value = -35530

# Row 94 (BSM)
$ws.Range("H94").Value = 2827.8096
$ws.Range("I94").Value = 1612.0667
$ws.Range("J94").Value = 5867.1665
$ws.Range("K94").Value = 1612.0667
$ws.Range("L94").Value = 5867.1665
$ws.Range("M94").Value = -1161.0667
$ws.Range("N94").Value = -6769.1665

# Row 105 (BSM)
$ws.Range("H105").Value = 24909.143
$ws.Range("I105").Value = 7899.6665
$ws.Range("K105").Value = 7899.6665
$ws.Range("M105").Value = -6152.6665

# Row 107 (BSM)
$ws.Range("H107").Value = 25756.666
$ws.Range("I107").Value = 28676.25
$ws.Range("K107").Value = 28676.25
$ws.Range("M107").Value = -26756.25

$ws = $wb.Worksheets.Item("CRP")
# Row 64 (CRP)
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

# Row 67 (CRP)
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

# Row 94 (CRP)
$ws.Range("H94").Value = 1481.375
$ws.Range("I94").Value = 2500
$ws.Range("J94").Value = 1335.8572
$ws.Range("K94").Value = 2500
$ws.Range("L94").Value = 1335.8572
$ws.Range("M94").Value = -2049
$ws.Range("N94").Value = -2237.8572

# Row 105 (CRP)
$ws.Range("H105").Value = 2638.5386
$ws.Range("I105").Value = 1644.5555
$ws.Range("K105").Value = 1644.5555
$ws.Range("M105").Value = 102.4445000000001

# Row 107 (CRP)
$ws.Range("H107").Value = 1894.2858
$ws.Range("I107").Value = 1894.2858
$ws.Range("K107").Value = 1894.2858
$ws.Range("M107").Value = 25.71419999999989

# Row 122 (CRP)
$ws.Range("H122").Value = 2862
$ws.Range("I122").Value = 2969.75
$ws.Range("K122").Value = 8909.25
$ws.Range("M122").Value = -6459.25

# Row 132 (CRP)
$ws.Range("H132").Value = 2826.3684
$ws.Range("I132").Value = 2537
$ws.Range("J132").Value = 3322.4285
$ws.Range("K132").Value = 7611
$ws.Range("L132").Value = 9967.2855
$ws.Range("M132").Value = -5081
$ws.Range("N132").Value = -15027.2855

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 1563.5
$ws.Range("I5").Value = 1227.125
$ws.Range("J5").Value = 1787.75
$ws.Range("K5").Value = 3681.375
$ws.Range("L5").Value = 5363.25
$ws.Range("M5").Value = -3569.375
$ws.Range("N5").Value = -5587.25

# Row 92 (CUL)
$ws.Range("H92").Value = 421.66666
$ws.Range("I92").Value = 385
$ws.Range("J92").Value = 495
$ws.Range("K92").Value = 1155
$ws.Range("L92").Value = 1485
$ws.Range("M92").Value = 93
$ws.Range("N92").Value = -3981

# Row 98 (CUL)
$ws.Range("H98").Value = 289.07693
$ws.Range("I98").Value = 250
$ws.Range("J98").Value = 292.33334
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 877.0000200000001
$ws.Range("M98").Value = 748
$ws.Range("N98").Value = -3873.00002

# Row 135 (CUL)
$ws.Range("H135").Value = 1563.5
$ws.Range("I135").Value = 1227.125
$ws.Range("J135").Value = 1787.75
$ws.Range("K135").Value = 11044.125
$ws.Range("L135").Value = 16089.75
$ws.Range("M135").Value = -8509.125
$ws.Range("N135").Value = -21159.75

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 84303.914
$ws.Range("I97").Value = 67414.336
$ws.Range("J97").Value = 112453.22
$ws.Range("K97").Value = 67414.336
$ws.Range("L97").Value = 112453.22
$ws.Range("M97").Value = -66918.336
$ws.Range("N97").Value = -113445.22

$ws = $wb.Worksheets.Item("LTW")
# Row 14 (LTW)
$ws.Range("H14").Value = 9959
$ws.Range("I14").Value = 9959
$ws.Range("K14").Value = 9959
$ws.Range("M14").Value = -9787

# Row 20 (LTW)
$ws.Range("H20").Value = 366675000
$ws.Range("J20").Value = 550000000
$ws.Range("L20").Value = 550000000
$ws.Range("N20").Value = -550000452

# Row 55 (LTW)
$ws.Range("H55").Value = 626.8684
$ws.Range("I55").Value = 390
$ws.Range("J55").Value = 781.34784
$ws.Range("K55").Value = 390
$ws.Range("L55").Value = 781.34784
$ws.Range("M55").Value = -217
$ws.Range("N55").Value = -1127.34784

# Row 100 (LTW)
$ws.Range("H100").Value = 4804.0586
$ws.Range("J100").Value = 4911.357
$ws.Range("L100").Value = 4911.357
$ws.Range("N100").Value = -5993.357

# Row 127 (LTW)
$ws.Range("H127").Value = 90546.625
$ws.Range("J127").Value = 90546.625
$ws.Range("L127").Value = 90546.625
$ws.Range("N127").Value = -100466.625

# Row 132 (LTW)
$ws.Range("H132").Value = 4437.091
$ws.Range("I132").Value = 5104.8335
$ws.Range("J132").Value = 3635.8
$ws.Range("K132").Value = 15314.5005
$ws.Range("L132").Value = 10907.4
$ws.Range("M132").Value = -12784.5005
$ws.Range("N132").Value = -15967.4

$ws = $wb.Worksheets.Item("WVR")
# Row 103 (WVR)
$ws.Range("H103").Value = 35432.168
$ws.Range("J103").Value = 35432.168
$ws.Range("L103").Value = 35432.168
$ws.Range("N103").Value = -37776.168

# Row 113 (WVR)
$ws.Range("H113").Value = 1268
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1268
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3804
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8144

# Row 132 (WVR)
$ws.Range("H132").Value = 2554.5518
$ws.Range("I132").Value = 2247.05
$ws.Range("J132").Value = 3237.889
$ws.Range("K132").Value = 6741.150000000001
$ws.Range("L132").Value = 9713.667000000001
$ws.Range("M132").Value = -4211.150000000001
$ws.Range("N132").Value = -14773.667

